# Add a new Job Posting row (Job_Id = 5) to Sheet1.
# Mirrors the commit "Add Job Posting with Job_Id=5": a brand-new row 6 is
# appended below the existing data (rows 1-5), extending the used range
# from A1:H5 to A1:H6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Dotnet Developer"
$ws.Range("C6").Value = "Hello"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
